$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "Partial cshtml..." text
# (paragraph 6 of the document) and build a Range that spans the whole
# paragraph, including its trailing paragraph mark, so that InsertXML
# *replaces* the paragraph (and the bookmark inside it) wholesale.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Partial cshtml scripts")) {
        $target = $p.Range
    }
}

$full = $d.Range($target.Start, $target.End)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = ""
$xml += "<w:p $ns><w:r><w:t>Partial cshtml scripts do not like a styles section. Maybe even a scripts section. They will need to be styled/scripted via tags in the same script</w:t></w:r></w:p>"
$xml += "<w:p $ns/>"
$xml += "<w:p $ns><w:r><w:t>It would seem that some things are easier in tag helpers and some things are easier with html helpers.</w:t></w:r></w:p>"
$xml += "<w:p $ns>"
$xml += '<w:r><w:t xml:space="preserve">For instance, </w:t></w:r>'
$xml += '<w:proofErr w:type="spellStart"/>'
$xml += '<w:r><w:t>its</w:t></w:r>'
$xml += '<w:proofErr w:type="spellEnd"/>'
$xml += '<w:r><w:t xml:space="preserve"> much easier and cleaner to write A tab strip</w:t></w:r>'
$xml += '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$xml += '<w:bookmarkEnd w:id="0"/>'
$xml += '<w:r><w:t xml:space="preserve"> with tag helpers. however  writing kendo buttons with html helpers is easier to pass parameters through for some JavaScript action </w:t></w:r>'
$xml += '<w:r><w:sym w:font="Wingdings" w:char="F04A"/></w:r>'
$xml += "</w:p>"
$xml += "<w:p $ns/>"

$full.InsertXML($xml)
